$d = $word.ActiveDocument

# Locate the existing Factorial Operation bullet text (end of the run) inside the
# "4. Feature Specifications" paragraph and append a new line-break + a second,
# shorter bullet describing the same Factorial Operation feature.
$rng = $d.Content
$found = $rng.Find.Execute(
    "Includes error handling for non-integer inputs, with appropriate feedback to users.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)  # wdCollapseEnd
    $newLine = [char]11 + "• Factorial Operation: Calculates the factorial of an integer using recursion or iteration and logs each step in history if applicable."
    $rng.InsertAfter($newLine)
}
